$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column widths (A -> ~11.78 chars, C -> ~10.78 chars)
$ws.Columns.Item(1).ColumnWidth = 11
$ws.Columns.Item(3).ColumnWidth = 10

# Row 4: StockName
$ws.Range("A4").Value = "StockName"

# Row 5: Price Table / OID / Quarter
# Shared-string pool order must be: Price Table, Quarter, OID
$ws.Range("C5").Value = "Price Table"
$ws.Range("E5").Value = "Quarter"
$ws.Range("D5").Value = "OID"

# Set the selection to L18 to match diff
$ws.Range("L18").Select() | Out-Null
